$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.151.11"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "2.488.91"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'321.78"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'108.97"
$ws.Range("E6").Value = "  +3.64%  "
$ws.Range("D7").Value = "'0.523"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("D10").Value = "'38.76"
$ws.Range("E10").Value = "  +7.03%  "
$ws.Range("D11").Value = "'0.0811"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "'18.23"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").Value = "'7.15"
$ws.Range("D15").Value = "2.880.43"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").Value = "2.487.17"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "47.086.25"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "'12.71"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").Value = "'6.61"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").Value = "0.0₃0936"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  +15.79%  "
$ws.Range("D23").Value = "'70.57"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'246.20"
$ws.Range("D25").Value = "'2.57"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'25.75"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").Value = "'2.28"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").Value = "'10.07"
$ws.Range("E29").Value = "  +2.84%  "
$ws.Range("E30").Value = "  +9.30%  "
$ws.Range("D31").Value = "'35.23"
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("D32").Value = "'49.91"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Value = "'20.08"
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("D34").Value = "'5.40"
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("D35").Value = "'0.0785"
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "'4.68"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("D39").Value = "'2.95"
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "'119.59"
$ws.Range("E42").Value = "  -3.26%  "
$ws.Range("D43").Value = "'21.39"
$ws.Range("E43").Value = "  +4.01%  "
$ws.Range("D44").Value = "'0.0295"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").Value = "1.984.47"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E46").Value = "  +1.95%  "
$ws.Range("E47").Value = "  -2.82%  "
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("D49").Value = "'9.06"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").Value = "'5.12"
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("D51").Value = "'57.22"
$ws.Range("E51").Value = "  +4.24%  "
